# Generate Report for Handoff
# This script updates the localization-status workbook so that the
# 6f7cbdfb-... file is reported as "Handed back: in sync with en-US"
# (row 2 on every sheet) and the 1fa6f70c-... file is reported as
# "Ready for handoff" (row 3 on every sheet), refreshing the related
# handoff/handback file names, datetimes and error detail message.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "6f7cbdfb-1204-463d-af87-57c34f87a784.md"
$ws1.Range("B2").Value = "e2e\6f7cbdfb-1204-463d-af87-57c34f87a784.md"

$ws1.Range("A3").Value = "1fa6f70c-acfd-4218-8e62-42b14564062a.md"
$ws1.Range("B3").Value = "e2e\1fa6f70c-acfd-4218-8e62-42b14564062a.md"

$ws1.Range("E3").Value = "Ready for handoff"
$ws1.Range("F3").Value = "Ready for handoff"
$ws1.Range("G3").Value = "2016-08-25 14:49:44"

# Hyperlink display text needs to be swapped too (the hyperlink
# addresses themselves stay pointed at the same targets).
$ws1_h1_addr = $ws1.Hyperlinks.Item(1).Address
$ws1_h2_addr = $ws1.Hyperlinks.Item(2).Address
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), $ws1_h1_addr, "", "", "e2e\6f7cbdfb-1204-463d-af87-57c34f87a784.md")
$ws1.Hyperlinks.Add($ws1.Range("B3"), $ws1_h2_addr, "", "", "e2e\1fa6f70c-acfd-4218-8e62-42b14564062a.md")

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "6f7cbdfb-1204-463d-af87-57c34f87a784.md"
$ws2.Range("G2").Value = "6f7cbdfb-1204-463d-af87-57c34f87a784.41c1c37f6b68f204dffca6cf14923904932fdd78.zh-cn.xlf"
$ws2.Range("I2").Value = "6f7cbdfb-1204-463d-af87-57c34f87a784.md"
$ws2.Range("J2").Value = "6f7cbdfb-1204-463d-af87-57c34f87a784.41c1c37f6b68f204dffca6cf14923904932fdd78.zh-cn.xlf"

$ws2.Range("A3").Value = "1fa6f70c-acfd-4218-8e62-42b14564062a.md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("G3").Value = "1fa6f70c-acfd-4218-8e62-42b14564062a.db3aae7004fcf9672d671966f88769de81053429.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-25 14:49:39"
$ws2.Range("I3").Value = "1fa6f70c-acfd-4218-8e62-42b14564062a.md"
$ws2.Range("J3").Value = "1fa6f70c-acfd-4218-8e62-42b14564062a.db3aae7004fcf9672d671966f88769de81053429.zh-cn.xlf"
$ws2.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aeba0416f269c6864473871c3d3f6a27df37ea43/e2e/1fa6f70c-acfd-4218-8e62-42b14564062a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb628b156bd1f6c3b74ee5648da0d02fa296a88a/e2e/1fa6f70c-acfd-4218-8e62-42b14564062a.md."

$ws2_h1_addr = $ws2.Hyperlinks.Item(1).Address
$ws2_h2_addr = $ws2.Hyperlinks.Item(2).Address
$ws2_h3_addr = $ws2.Hyperlinks.Item(3).Address
$ws2_h4_addr = $ws2.Hyperlinks.Item(4).Address
$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $ws2_h1_addr, "", "", "6f7cbdfb-1204-463d-af87-57c34f87a784.md")
$ws2.Hyperlinks.Add($ws2.Range("I2"), $ws2_h2_addr, "", "", "6f7cbdfb-1204-463d-af87-57c34f87a784.md")
$ws2.Hyperlinks.Add($ws2.Range("A3"), $ws2_h3_addr, "", "", "1fa6f70c-acfd-4218-8e62-42b14564062a.md")
$ws2.Hyperlinks.Add($ws2.Range("I3"), $ws2_h4_addr, "", "", "1fa6f70c-acfd-4218-8e62-42b14564062a.md")

# Error Detail column grew a long message; widen it like the other
# "long text" columns in the sheet (A, G, I, J already use width 40).
$ws2.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "6f7cbdfb-1204-463d-af87-57c34f87a784.md"
$ws3.Range("G2").Value = "6f7cbdfb-1204-463d-af87-57c34f87a784.41c1c37f6b68f204dffca6cf14923904932fdd78.de-de.xlf"
$ws3.Range("I2").Value = "6f7cbdfb-1204-463d-af87-57c34f87a784.md"
$ws3.Range("J2").Value = "6f7cbdfb-1204-463d-af87-57c34f87a784.41c1c37f6b68f204dffca6cf14923904932fdd78.de-de.xlf"

$ws3.Range("A3").Value = "1fa6f70c-acfd-4218-8e62-42b14564062a.md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("G3").Value = "1fa6f70c-acfd-4218-8e62-42b14564062a.db3aae7004fcf9672d671966f88769de81053429.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-25 14:49:44"
$ws3.Range("I3").Value = "1fa6f70c-acfd-4218-8e62-42b14564062a.md"
$ws3.Range("J3").Value = "1fa6f70c-acfd-4218-8e62-42b14564062a.db3aae7004fcf9672d671966f88769de81053429.de-de.xlf"
$ws3.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aeba0416f269c6864473871c3d3f6a27df37ea43/e2e/1fa6f70c-acfd-4218-8e62-42b14564062a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb628b156bd1f6c3b74ee5648da0d02fa296a88a/e2e/1fa6f70c-acfd-4218-8e62-42b14564062a.md."

$ws3_h1_addr = $ws3.Hyperlinks.Item(1).Address
$ws3_h2_addr = $ws3.Hyperlinks.Item(2).Address
$ws3_h3_addr = $ws3.Hyperlinks.Item(3).Address
$ws3_h4_addr = $ws3.Hyperlinks.Item(4).Address
$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $ws3_h1_addr, "", "", "6f7cbdfb-1204-463d-af87-57c34f87a784.md")
$ws3.Hyperlinks.Add($ws3.Range("I2"), $ws3_h2_addr, "", "", "6f7cbdfb-1204-463d-af87-57c34f87a784.md")
$ws3.Hyperlinks.Add($ws3.Range("A3"), $ws3_h3_addr, "", "", "1fa6f70c-acfd-4218-8e62-42b14564062a.md")
$ws3.Hyperlinks.Add($ws3.Range("I3"), $ws3_h4_addr, "", "", "1fa6f70c-acfd-4218-8e62-42b14564062a.md")

$ws3.Columns.Item(16).ColumnWidth = 39.17
